$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows 56-105 (A:D pattern continues the boundary-condition table; rows
#     56-59 additionally carry node-count values/formula in column H) ---

$ws.Cells.Item(56,1).Value = 2
$ws.Cells.Item(56,2).Value = 0
$ws.Cells.Item(56,3).Value = 1
$ws.Cells.Item(56,4).Value = 0
$ws.Cells.Item(56,8).Value = 50

$ws.Cells.Item(57,1).Value = 3
$ws.Cells.Item(57,2).Value = 0
$ws.Cells.Item(57,3).Value = 1
$ws.Cells.Item(57,4).Value = 0
$ws.Cells.Item(57,8).Value = 1

$ws.Cells.Item(58,1).Value = 4
$ws.Cells.Item(58,2).Value = 0
$ws.Cells.Item(58,3).Value = 1
$ws.Cells.Item(58,4).Value = 0
$ws.Cells.Item(58,8).Value = 3

$ws.Cells.Item(59,1).Value = 5
$ws.Cells.Item(59,2).Value = 0
$ws.Cells.Item(59,3).Value = 1
$ws.Cells.Item(59,4).Value = 0
$ws.Cells.Item(59,8).Formula = "=SUM(H56:H58)"

# Rows 60-105: A counts up from 6 to 51, B/D are 0, C is 1
$a = 6
for ($r = 60; $r -le 105; $r++) {
    $ws.Cells.Item($r,1).Value = $a
    $ws.Cells.Item($r,2).Value = 0
    $ws.Cells.Item($r,3).Value = 1
    $ws.Cells.Item($r,4).Value = 0
    $a++
}

# Recalculate so the new SUM formula's cached value is correct
$excel.Calculate()

# --- View state: update selection / scroll position / zoom to match the edited sheet ---
$ws.Activate() | Out-Null
$ws.Range("K65").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 100
